# "Generate Report for handback" -- refresh the handback status report.
#
# For each localized-language sheet (zh-cn, de-de) in the workbook:
#   - Status (col B) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every file row.
#   - Two new columns get populated for each file row:
#       E: Latest Target File    (same file as the source .md, hyperlinked)
#       F: Latest Handback File  (same file as the handoff .xlf, hyperlinked)
#   - Latest Handback DateTime (col G) is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# Per-language sheet data: handoff-time (existing col D, untouched) is not
# needed here -- only the new handback timestamp and the hyperlink targets
# for the two data rows (186e4dd2... and 1edf6f82...) differ per sheet.
$sheetsInfo = @(
    @{
        Name = "zh-cn"
        HandbackTime = "2016-01-21 03:01:31"
        Rows = @(
            @{
                Row = 2
                MdName = "186e4dd2-2f28-406b-a12e-58c65bddabed.md"
                MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/abcb6eab6bd92f420e7f0d5d22a9a44aef601707/e2e/186e4dd2-2f28-406b-a12e-58c65bddabed.md"
                XlfName = "186e4dd2-2f28-406b-a12e-58c65bddabed.5c63508d08cb0d5141fdc1b663117b1cb059d217.zh-cn.xlf"
                XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43a2cc71e5b6bb0f75fb918a1f0e84ab7dcae5f4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/186e4dd2-2f28-406b-a12e-58c65bddabed.5c63508d08cb0d5141fdc1b663117b1cb059d217.zh-cn.xlf"
            },
            @{
                Row = 3
                MdName = "1edf6f82-823e-4699-8a91-4598dcfa9b3c.md"
                MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/abcb6eab6bd92f420e7f0d5d22a9a44aef601707/e2e/1edf6f82-823e-4699-8a91-4598dcfa9b3c.md"
                XlfName = "1edf6f82-823e-4699-8a91-4598dcfa9b3c.c58b9a3b30f234b06e1b26f163084c99920a4b85.zh-cn.xlf"
                XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43a2cc71e5b6bb0f75fb918a1f0e84ab7dcae5f4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/1edf6f82-823e-4699-8a91-4598dcfa9b3c.c58b9a3b30f234b06e1b26f163084c99920a4b85.zh-cn.xlf"
            }
        )
    },
    @{
        Name = "de-de"
        HandbackTime = "2016-01-21 03:01:52"
        Rows = @(
            @{
                Row = 2
                MdName = "186e4dd2-2f28-406b-a12e-58c65bddabed.md"
                MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/abcb6eab6bd92f420e7f0d5d22a9a44aef601707/e2e/186e4dd2-2f28-406b-a12e-58c65bddabed.md"
                XlfName = "186e4dd2-2f28-406b-a12e-58c65bddabed.5c63508d08cb0d5141fdc1b663117b1cb059d217.de-de.xlf"
                XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/23ba0560160f8b6845cc1d99aabe589800f4bd09/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/186e4dd2-2f28-406b-a12e-58c65bddabed.5c63508d08cb0d5141fdc1b663117b1cb059d217.de-de.xlf"
            },
            @{
                Row = 3
                MdName = "1edf6f82-823e-4699-8a91-4598dcfa9b3c.md"
                MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/abcb6eab6bd92f420e7f0d5d22a9a44aef601707/e2e/1edf6f82-823e-4699-8a91-4598dcfa9b3c.md"
                XlfName = "1edf6f82-823e-4699-8a91-4598dcfa9b3c.c58b9a3b30f234b06e1b26f163084c99920a4b85.de-de.xlf"
                XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/23ba0560160f8b6845cc1d99aabe589800f4bd09/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/1edf6f82-823e-4699-8a91-4598dcfa9b3c.c58b9a3b30f234b06e1b26f163084c99920a4b85.de-de.xlf"
            }
        )
    }
)

foreach ($sheetInfo in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    foreach ($rowInfo in $sheetInfo.Rows) {
        $r = $rowInfo.Row

        # Status: handed back, in sync with en-US.
        $ws.Cells.Item($r, 2).Value = $status

        # New "Latest Target File" column (E) -- hyperlink to the source .md,
        # same as column A for this row.
        $eCell = $ws.Cells.Item($r, 5)
        $ws.Hyperlinks.Add($eCell, $rowInfo.MdUrl, "", "", $rowInfo.MdName) | Out-Null

        # New "Latest Handback File" column (F) -- hyperlink to the handoff
        # .xlf, same as column C for this row.
        $fCell = $ws.Cells.Item($r, 6)
        $ws.Hyperlinks.Add($fCell, $rowInfo.XlfUrl, "", "", $rowInfo.XlfName) | Out-Null

        # Latest Handback DateTime (G) -- stamp with the handback time.
        $ws.Cells.Item($r, 7).Value = $sheetInfo.HandbackTime
    }
}
